$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 is a date-shaped string ("03/08/2025") that must stay TEXT, not be
# auto-converted to a date serial number. Force text entry via NumberFormat,
# then reset the cell style back to Normal so no stray number format lingers.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "03/08/2025"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "Celtic"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "St. Mirren"
$ws.Range("F2").Value = "W"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1.73
$ws.Range("L2").Value = 0.27
$ws.Range("M2").Value = 14
$ws.Range("N2").Value = 4
$ws.Range("O2").Value = 6
$ws.Range("P2").Value = 2
